$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: new METRO entry 'PAL1'
$ws.Range("A5").Value = "PAL1"
$ws.Range("B5").Value = "METRO"
$ws.Range("C5").Value = "ME9-D7-PAL1"
$ws.Range("D5").Value = "172.31.250.96"
$ws.Range("E5").Value = "username_nms"
$ws.Range("F5").Value = "password_nms"

# Row 8: new Switch entry 'DGL007'
$ws.Range("A8").Value = "DGL007"
$ws.Range("B8").Value = "L2SW FH S5800"
$ws.Range("C8").Value = "L2SW-D7-DGL007"
$ws.Range("D8").Value = "10.199.162.2"
$ws.Range("E8").Value = "admin"
$ws.Range("F8").Value = 12345

# Row 12: DGL034 (F entered as text "12345")
$ws.Range("A12").Value = "DGL034"
$ws.Range("B12").Value = "L2SW FH S5800"
$ws.Range("C12").Value = "L2SW-D7-DGL034"
$ws.Range("D12").Value = "10.199.162.78"
$ws.Range("E12").Value = "admin"
$ws.Range("F12").Value = "12345"

# Row 4: rename PGI -> PRG
$ws.Range("A4").Value = "PRG"
$ws.Range("B4").Value = "METRO"
$ws.Range("C4").Value = "ME-D7-PRG"
$ws.Range("D4").Value = "172.32.250.101"
$ws.Range("E4").Value = "username_nms"
$ws.Range("F4").Value = "password_nms"

$ws.Range("A6").Value = "STG"
$ws.Range("B6").Value = "METRO"
$ws.Range("C6").Value = "ME-D7-STG"
$ws.Range("D6").Value = "172.32.250.102"
$ws.Range("E6").Value = "username_nms"
$ws.Range("F6").Value = "password_nms"

$ws.Range("A7").Value = "DGL006"
$ws.Range("B7").Value = "L2SW RAISECOM"
$ws.Range("C7").Value = "L2SW-D7-DGL006"
$ws.Range("D7").Value = "172.25.88.18"
$ws.Range("E7").Value = "raisecom"
$ws.Range("F7").Value = "raisecom"

$ws.Range("A9").Value = "DGL105"
$ws.Range("B9").Value = "L2SW RAISECOM"
$ws.Range("C9").Value = "L2SW-D7-DGL105"
$ws.Range("D9").Value = "10.199.162.85"
$ws.Range("E9").Value = "raisecom"
$ws.Range("F9").Value = "raisecom"

$ws.Range("A10").Value = "DGL129"
$ws.Range("B10").Value = "L2SW RAISECOM"
$ws.Range("C10").Value = "L2SW-D7-DGL129"
$ws.Range("D10").Value = "172.25.88.19"
$ws.Range("E10").Value = "raisecom"
$ws.Range("F10").Value = "raisecom"

$ws.Range("A11").Value = "DGL030"
$ws.Range("B11").Value = "L2SW FH S5800"
$ws.Range("C11").Value = "L2SW-D7-DGL030"
$ws.Range("D11").Value = "10.199.162.25"
$ws.Range("E11").Value = "admin"
$ws.Range("F11").Value = 12345

$ws.Range("A13").Value = "DGL175"
$ws.Range("B13").Value = "L2SW FH S5800"
$ws.Range("C13").Value = "L2SW-D7-DGL175"
$ws.Range("D13").Value = "10.199.162.66"
$ws.Range("E13").Value = "admin"
$ws.Range("F13").Value = 12345

$ws.Range("A14").Value = "DGL079"
$ws.Range("B14").Value = "L2SW FH S5800"
$ws.Range("C14").Value = "L2SW-D7-DGL079"
$ws.Range("D14").Value = "10.199.162.26"
$ws.Range("E14").Value = "admin"
$ws.Range("F14").Value = 12345

$ws.Range("A15").Value = "PGI003"
$ws.Range("B15").Value = "L2SW FH CITRANS"
$ws.Range("C15").Value = "SW-D7-TSEL-PGI003-10G-4"
$ws.Range("D15").Value = "172.25.88.7"
$ws.Range("E15").Value = "username_nms"
$ws.Range("F15").Value = "password_nms"

$ws.Range("A16").Value = "PGI004"
$ws.Range("B16").Value = "L2SW FH S6800"
$ws.Range("C16").Value = "SW-D7-PGI004"
$ws.Range("D16").Value = "10.198.2.130"
$ws.Range("E16").Value = "admin"
$ws.Range("F16").Value = "Admin12345"

$ws.Range("A17").Value = "PGI063"
$ws.Range("B17").Value = "L2SW RAISECOM"
$ws.Range("C17").Value = "L2SW-D7-PGI063"
$ws.Range("D17").Value = "10.199.162.49"
$ws.Range("E17").Value = "raisecom"
$ws.Range("F17").Value = "raisecom"

$ws.Range("A18").Value = "TLI005"
$ws.Range("B18").Value = "L2SW RAISECOM"
$ws.Range("C18").Value = "L2SW-D7-TLI005"
$ws.Range("D18").Value = "10.199.162.76"
$ws.Range("E18").Value = "raisecom"
$ws.Range("F18").Value = "raisecom"

$ws.Range("A19").Value = "TLI008"
$ws.Range("B19").Value = "L2SW FH S5800"
$ws.Range("C19").Value = "SW-D7-TSEL-TLI008"
$ws.Range("D19").Value = "10.199.162.30"
$ws.Range("E19").Value = "admin"
$ws.Range("F19").Value = 12345

$ws.Range("A20").Value = "TLI041"
$ws.Range("B20").Value = "L2SW RAISECOM"
$ws.Range("C20").Value = "SW-D7-TSEL-TLI041"
$ws.Range("D20").Value = "10.199.162.32"
$ws.Range("E20").Value = "raisecom"
$ws.Range("F20").Value = "raisecom"

$ws.Range("I9").Select()
